$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Structural change: the "Odd_CS_3-3_HT" column (formerly at BC) is moved to
# become the first of the "_HT" correct-score odds columns (now at AW),
# shifting the intervening columns (AW:BB -> AX:BC) one position to the right.
$ws.Columns.Item(55).Cut()
$ws.Columns.Item(49).Insert()

# Updated odds values throughout the sheet
$ws.Range("G2").Value = 1.95
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 2.6
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 3.75
$ws.Range("S2").Value = 1.25
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 11
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 19
$ws.Range("AA2").Value = 15
$ws.Range("AD2").Value = 7
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 12
$ws.Range("AL2").Value = 23
$ws.Range("AM2").Value = 26
$ws.Range("AO2").Value = 10
$ws.Range("AQ2").Value = 34
$ws.Range("AU2").Value = 7
$ws.Range("AW2").Value = 351
$ws.Range("AX2").Value = 5.5
$ws.Range("AY2").Value = 17
$ws.Range("AZ2").Value = 21
$ws.Range("W3").Value = 9.25
$ws.Range("X3").Value = 13
$ws.Range("AC3").Value = 11.25
$ws.Range("AD3").Value = 5.5
$ws.Range("AH3").Value = 8.75
$ws.Range("AM3").Value = 18.5
$ws.Range("AN3").Value = 4.7
$ws.Range("AP3").Value = 17.5
$ws.Range("AR3").Value = 70
$ws.Range("AT3").Value = 2.82
$ws.Range("AU3").Value = 6.1
$ws.Range("AZ3").Value = 19
$ws.Range("BB3").Value = 80
$ws.Range("H4").Value = 3.8
$ws.Range("I4").Value = 3.9
$ws.Range("J4").Value = 2.2
$ws.Range("K4").Value = 2.27
$ws.Range("L4").Value = 4.2
$ws.Range("R4").Value = 2.05
$ws.Range("W4").Value = 7.2
$ws.Range("X4").Value = 7.6
$ws.Range("Y4").Value = 7.1
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 6.7
$ws.Range("AE4").Value = 11.75
$ws.Range("AH4").Value = 11.5
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 11.25
$ws.Range("AK4").Value = 45
$ws.Range("AL4").Value = 27
$ws.Range("AO4").Value = 8
$ws.Range("AQ4").Value = 26
$ws.Range("AT4").Value = 3.05
$ws.Range("AU4").Value = 7.2
$ws.Range("AX4").Value = 5.9
$ws.Range("AY4").Value = 21
$ws.Range("BA4").Value = 110
$ws.Range("BB4").Value = 120
$ws.Range("G5").Value = 2.35
$ws.Range("I5").Value = 2.63
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("Y5").Value = 9.5
$ws.Range("AA5").Value = 17
$ws.Range("AJ5").Value = 11
$ws.Range("G6").Value = 1.7
$ws.Range("H6").Value = 3.65
$ws.Range("I6").Value = 4.45
$ws.Range("J6").Value = 2.22
$ws.Range("K6").Value = 2.18
$ws.Range("L6").Value = 4.75
$ws.Range("N6").Value = 11
$ws.Range("Q6").Value = 1.72
$ws.Range("R6").Value = 1.88
$ws.Range("U6").Value = 1.72
$ws.Range("V6").Value = 1.9
$ws.Range("W6").Value = 7.4
$ws.Range("X6").Value = 8.25
$ws.Range("Z6").Value = 13.5
$ws.Range("AA6").Value = 13
$ws.Range("AB6").Value = 24
$ws.Range("AC6").Value = 11
$ws.Range("AD6").Value = 7.2
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 65
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 13
$ws.Range("AI6").Value = 26
$ws.Range("AJ6").Value = 14.5
$ws.Range("AL6").Value = 45
$ws.Range("AN6").Value = 3.55
$ws.Range("AO6").Value = 8.25
$ws.Range("AP6").Value = 16.5
$ws.Range("AQ6").Value = 27
$ws.Range("AT6").Value = 2.77
$ws.Range("AU6").Value = 7.3
$ws.Range("AX6").Value = 6.2
$ws.Range("AY6").Value = 26
$ws.Range("G7").Value = 1.65
$ws.Range("I7").Value = 5.2
$ws.Range("P7").Value = 2.77
$ws.Range("T7").Value = 2.55
$ws.Range("W7").Value = 5.9
$ws.Range("X7").Value = 7
$ws.Range("AA7").Value = 14
$ws.Range("AD7").Value = 6.8
$ws.Range("AI7").Value = 30
$ws.Range("AK7").Value = 100
$ws.Range("AM7").Value = 65
$ws.Range("AN7").Value = 3.4
$ws.Range("AT7").Value = 2.52
$ws.Range("AV7").Value = 75
$ws.Range("AX7").Value = 6.7
$ws.Range("BB7").Value = 250
$ws.Range("BC7").Value = 500
$ws.Range("N8").Value = 12.8
$ws.Range("S8").Value = 1.28
$ws.Range("T8").Value = 3.34
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 3.9
$ws.Range("K9").Value = 2.05
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.7
$ws.Range("U9").Value = 1.95
$ws.Range("V9").Value = 1.8
$ws.Range("W9").Value = 6.5
$ws.Range("X9").Value = 9
$ws.Range("AC9").Value = 8.5
$ws.Range("AD9").Value = 6
$ws.Range("AG9").Value = 351
$ws.Range("AL9").Value = 34
$ws.Range("AO9").Value = 12
$ws.Range("AR9").Value = 67
$ws.Range("AS9").Value = 201
$ws.Range("AU9").Value = 8.5
$ws.Range("AV9").Value = 67
$ws.Range("AZ9").Value = 34
$ws.Range("BC9").Value = 251
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 3.55
$ws.Range("I11").Value = 1.62
$ws.Range("J11").Value = 5.4
$ws.Range("K11").Value = 2.15
$ws.Range("L11").Value = 2.18
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 6.7
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 2.92
$ws.Range("Q11").Value = 2.05
$ws.Range("R11").Value = 1.7
$ws.Range("T11").Value = 2.67
$ws.Range("U11").Value = 2.02
$ws.Range("V11").Value = 1.7
$ws.Range("W11").Value = 11.75
$ws.Range("X11").Value = 28
$ws.Range("Y11").Value = 17
$ws.Range("Z11").Value = 100
$ws.Range("AA11").Value = 60
$ws.Range("AB11").Value = 70
$ws.Range("AC11").Value = 6.7
$ws.Range("AD11").Value = 7
$ws.Range("AE11").Value = 19
$ws.Range("AF11").Value = 110
$ws.Range("AI11").Value = 6.9
$ws.Range("AJ11").Value = 8.25
$ws.Range("AK11").Value = 11.75
$ws.Range("AL11").Value = 14
$ws.Range("AN11").Value = 6.7
$ws.Range("AO11").Value = 32
$ws.Range("AP11").Value = 40
$ws.Range("AQ11").Value = 200
$ws.Range("AT11").Value = 2.67
$ws.Range("AU11").Value = 8
$ws.Range("AX11").Value = 3.4
$ws.Range("AY11").Value = 7.8
$ws.Range("AZ11").Value = 18.5
$ws.Range("BA11").Value = 26
$ws.Range("BB11").Value = 60
